# Update countries & provincias Spain
# Applies the COVID data refresh captured by the source diff:
#  - four pairs of countries swap rank (their row's stats are replaced
#    with the newer figures, causing the two neighbouring rows to trade
#    country names as the table stays sorted by "Casos totales" desc.)
#  - six more rows get refreshed case/death figures with no reordering
#  - the "last updated" timestamp banner in A1 is bumped

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes)

    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Rows whose figures refresh but keep their current country (row order unchanged)
Set-Row 7   "Rusia"      1091186 5905 901207 170784 0 134 19195
Set-Row 24  "Filipinas"  279526  3257 208790 65906  0 47  4830
Set-Row 57  "Singapur"   57543   11   57039  477    0 0   27
Set-Row 67  "Afganistan" 38883   11   32576  4870   0 1   1437
Set-Row 144 "Estonia"    2814    36   2357   393    0 0   64
Set-Row 160 "Letonia"    1498    4    1248   214    0 0   36

# Rows that swap rank with their neighbour: row gets new leader's data,
# country name updates to match, and the displaced country's old figures
# move down into the following row.
Set-Row 89  "Croacia"        14513 234 12169 2100  0 6 244
Set-Row 90  "Grecia"         14400 0   3804  10271 0 0 325

Set-Row 110 "Eslovaquia"     6256  235 3390  2827  0 0 39
Set-Row 111 "Mozambique"     6161  0   3393  2729  0 0 39

Set-Row 131 "Lituania"       3565  61  2181  1297  0 0 87
Set-Row 132 "Mayotte"        3541  0   2964  537   0 0 40

Set-Row 214 "Islas Malvinas" 13    0   13    0     0 0 0
Set-Row 215 "Montserrat"     13    0   12    0     0 0 1

# Refresh the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 10:13"
